# Automatische test-sync: 2025-08-03 18:40:50
# Adds a new incoming-mail log row (#13) to the "Logs" sheet and refreshes
# the "Dashboard" summary sheet to reflect the updated category counts.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Append the new log entry as row 41 on the "Logs" sheet
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")
$newRow = 41

$logs.Cells.Item($newRow, 1).Value  = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value  = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Cells.Item($newRow, 4).Value  = "Documentatie / Datasheets"
$logs.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar documentatie@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value  = "2025-08-03 18:40:14"
$logs.Cells.Item($newRow, 7).Value  = "Ja"
$logs.Cells.Item($newRow, 8).Value  = "Ja"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# ---------------------------------------------------------------------
# 2. Extend the conditional-formatting ranges from row 40 to row 41
#    so the newly added row gets the same colour rules as the rest.
# ---------------------------------------------------------------------
$ccols = @("D", "G", "H", "I", "J")
foreach ($col in $ccols) {
    $oldRange = $logs.Range($col + "2:" + $col + "40")
    $newRange = $logs.Range($col + "2:" + $col + "41")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Refresh the "Dashboard" sheet: the "Documentatie / Datasheets"
#    category now ties with "Retour / Terugbetaling" at 2 occurrences,
#    so it moves up to row 6 and "Retour / Terugbetaling" drops to row 7
#    with its count updated to 2.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Documentatie / Datasheets"
$dash.Range("B6").Value = 2
$dash.Range("A7").Value = "Retour / Terugbetaling"
$dash.Range("B7").Value = 2
